# Update countries & provincias Spain
#
# This script reproduces the daily data-refresh edit of the COVID paises.xlsx
# workbook:
#   - Chile's numbers grew enough that it now ranks above Irlanda/Pakistan,
#     so those three rows (26-28) shift down and Chile's fresh numbers land
#     in row 26.
#   - Egipto's numbers grew enough that it now ranks above Australia, so
#     those two rows (51-52) swap.
#   - A handful of other country rows simply get refreshed numbers.
#   - The "last updated" timestamp banner moves from 17:03 to 17:33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Header timestamp -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 17:33"

# --- Estados Unidos (row 4) refreshed totals ------------------------------
Set-Row 4 "Estados Unidos" 1216150 3315 188075 957914 16055 240 70161

# --- Chile overtakes Irlanda and Pakistan (rows 26-28) --------------------
Set-Row 26 "Chile" 22016 1373 10710 11031 470 5 275
Set-Row 27 "Irlanda" 21772 0 13386 7067 93 0 1319
Set-Row 28 "Pakistan" 21501 560 5782 15233 111 10 486

# --- Republica Dominicana (row 45) refreshed totals -----------------------
Set-Row 45 "Republica Dominicana" 8480 245 1905 6221 144 8 354

# --- Egipto overtakes Australia (rows 51-52) -------------------------------
Set-Row 51 "Egipto" 7201 388 1730 5019 0 16 452
Set-Row 52 "Australia" 6849 24 5889 864 27 1 96

# --- Argelia (row 58) refreshed totals -------------------------------------
Set-Row 58 "Argelia" 4838 190 2067 2301 22 5 470

# --- Moldavia (row 59) refreshed totals -------------------------------------
Set-Row 59 "Moldavia" 4363 115 1544 2683 237 4 136

# --- Republica de Chipre (row 95) refreshed totals --------------------------
Set-Row 95 "Republica de Chipre" 878 4 296 567 15 0 15
